# Updates cryptos list values (price/volume) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.635.79'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.585.92'
$ws.Range('E3').Value = '  -2.75%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'210.94"
$ws.Range('E5').Value = '  -2.39%  '
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -2.63%  '
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').Value = "'19.62"
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').Value = '1.806.81'
$ws.Range('E12').Value = '  -2.82%  '
$ws.Range('D13').Value = '1.590.73'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').Value = "'0.528"
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').Value = "'64.63"
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '26.624.61'
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'208.12"
$ws.Range('E19').Value = '  -3.79%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = "'1.00"
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = "'6.76"
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('E23').Value = '  -4.98%  '
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').Value = "'146.81"
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'7.41"
$ws.Range('E27').Value = '  +1.67%  '
$ws.Range('D28').Value = "'0.114"
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').Value = "'15.32"
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  -1.91%  '
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('D33').Value = "'0.683"
$ws.Range('E33').Value = '  +26.14%  '
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').Value = '1.331.94'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.44"
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = "'1.51"
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').Value = "'0.826"
$ws.Range('E39').Value = '  -2.85%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('D42').Value = "'0.784"
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('E43').Value = '  -3.75%  '
$ws.Range('D44').Value = "'63.56"
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').Value = '1.720.72'
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').Value = "'89.76"
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').Value = "'0.828"
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').Value = "'7.46"
$ws.Range('E51').Value = '  -1.31%  '
